$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the date serial 45172 (2023-09-03) for every
# data row (rows 2 through 238). Update it to 45175 (2023-09-06) for all of
# them in one shot using the underlying date serial value.
$ws.Range("C2:C238").Value = 45175
